# feat: xls com oportunidades consolidada + analises
#
# Adds a "Probabilidade" column (E) computed as 1/(Odd+Margem) on the
# "Oportunidades" sheet, shifting the existing "Resultado" column to F and
# filling in the actual results. Formats the new column as a percentage and
# applies a style to H7 (stray formatted cell left over from formatting the
# new range). Also restores window geometry / file-version style metadata
# that Excel rewrites on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oportunidades")
$ws.Activate()

# Insert a new column before the current "Resultado" column (E) so the
# layout becomes: Jogo | Market | Margem | Odd | Probabilidade | Resultado
$ws.Range("E1").EntireColumn.Insert()

# Header for the new column
$ws.Range("E1").Value = "Probabilidade"
$ws.Range("E1").Style = $ws.Range("F1").Style

# Formulas for the probability column: 1 / (Odd + Margem)
$ws.Range("E2:E8").Formula = "=1/(D2+C2)"

# Apply a percentage number format (style "Porcentagem") to the new column
$percentStyle = $wb.Styles.Add("Porcentagem")
$percentStyle.NumberFormat = "0%"
$ws.Range("E2:E8").Style = "Porcentagem"

# Fill in the actual results (previously all placeholder zeros) in the
# now-shifted "Resultado" column (F)
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 0.3
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -1

# Column F width (bestFit)
$ws.Range("F1").EntireColumn.AutoFit()

# Leftover styled (but empty) cell at H7, underlined font, matching the diff
$ws.Range("H7").Font.Underline = $true
$ws.Range("H7").Value = $null

# Selection as left after editing
$ws.Range("E2").Select()

# Page setup (paper size / orientation) as captured by the diff
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
